$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'23.531.83"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").Value = "'1.650.46"
$ws.Range("E3").Value = "  +0.68%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("E5").Value = "  +0.22%  "
$ws.Range("D6").Value = "'300.63"
$ws.Range("E6").Value = "  -1.03%  "
$ws.Range("D7").Value = "'0.3786"
$ws.Range("E7").Value = "  +0.64%  "
$ws.Range("D8").Value = "'0.3577"
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").Value = "'50.74"
$ws.Range("E9").Value = "  -1.45%  "
$ws.Range("D10").Value = "'0.08115"
$ws.Range("D11").Value = "'1.226"
$ws.Range("E11").Value = "  -0.34%  "
$ws.Range("D12").Value = "'1.002"
$ws.Range("E12").Value = "  +0.32%  "
$ws.Range("D13").Value = "'22.13"
$ws.Range("E13").Value = "  -0.84%  "
$ws.Range("D14").Value = "'6.425"
$ws.Range("E14").Value = "  -1.38%  "
$ws.Range("D15").Value = "'7.422"
$ws.Range("E15").Value = "  +1.23%  "
$ws.Range("D16").Value = "'0.00001205"
$ws.Range("E16").Value = "  -1.69%  "
$ws.Range("D17").Value = "'1.659.65"
$ws.Range("E17").Value = "  +1.49%  "
$ws.Range("D18").Value = "'97.11"
$ws.Range("E19").Value = "  +0.56%  "
$ws.Range("D20").Value = "'6.788"
$ws.Range("E20").Value = "  +0.76%  "
$ws.Range("D21").Value = "'17.48"
$ws.Range("E21").Value = "  +0.13%  "
$ws.Range("E22").Value = "  +0.21%  "
$ws.Range("D23").Value = "'12.64"
$ws.Range("E23").Value = "  +1.07%  "
$ws.Range("D24").Value = "'23.557.86"
$ws.Range("E24").Value = "  +0.44%  "
$ws.Range("D25").Value = "'2.491"
$ws.Range("E25").Value = "  -1.11%  "
$ws.Range("D26").Value = "'2.940"
$ws.Range("E26").Value = "  -5.61%  "
$ws.Range("D27").Value = "'20.98"
$ws.Range("E27").Value = "  -0.91%  "
$ws.Range("D28").Value = "'153.09"
$ws.Range("E28").Value = "  +0.35%  "
$ws.Range("D29").Value = "'5.236"
$ws.Range("E29").Value = "  +0.70%  "
$ws.Range("D30").Value = "'133.30"
$ws.Range("E30").Value = "  -0.54%  "
$ws.Range("D31").Value = "'1.833.47"
$ws.Range("E31").Value = "  +0.82%  "
$ws.Range("D32").Value = "'6.999"
$ws.Range("E32").Value = "  +3.82%  "
$ws.Range("D33").Value = "'2.146"
$ws.Range("E33").Value = "  +5.34%  "
$ws.Range("D34").Value = "'11.92"
$ws.Range("E34").Value = "  +3.81%  "
$ws.Range("E35").Value = "  -4.79%  "
$ws.Range("D36").Value = "'0.02743"
$ws.Range("E36").Value = "  -0.73%  "
$ws.Range("D37").Value = "'0.08737"
$ws.Range("E37").Value = "  -0.34%  "
$ws.Range("D38").Value = "'0.2457"
$ws.Range("E38").Value = "  -1.32%  "
$ws.Range("D39").Value = "'5.999"
$ws.Range("E39").Value = "  +0.07%  "
$ws.Range("D40").Value = "'13.18"
$ws.Range("E40").Value = "  +4.51%  "
$ws.Range("E41").Value = "  -1.37%  "
$ws.Range("D42").Value = "'0.6936"
$ws.Range("E42").Value = "  -0.77%  "
$ws.Range("D43").Value = "'1.324"
$ws.Range("E43").Value = "  +0.16%  "
$ws.Range("E44").Value = "  +0.67%  "
$ws.Range("E45").Value = "  +0.16%  "
$ws.Range("E46").Value = "  +0.20%  "
$ws.Range("D47").Value = "'2.270"
$ws.Range("E47").Value = "  -2.26%  "
$ws.Range("D48").Value = "'3.933"
$ws.Range("E48").Value = "  -0.54%  "
$ws.Range("D49").Value = "'0.07824"
$ws.Range("E49").Value = "  -1.38%  "
$ws.Range("D50").Value = "'128.43"
$ws.Range("E50").Value = "  +0.70%  "
$ws.Range("D51").Value = "'1.174"
$ws.Range("E51").Value = "  -0.41%  "
